# Apply crypto price/volume updates per the commit diff.
# Numeric-looking text values get an apostrophe (quote) prefix so Excel
# keeps them as text (matching the original inlineStr cell type) instead
# of silently converting them to numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.976.15"
$ws.Range("E2").Value = "  -2.92%  "
$ws.Range("D3").Value = "3.837.98"
$ws.Range("E3").Value = "  -2.55%  "
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("D5").Value = "'599.92"
$ws.Range("E5").Value = "  -1.91%  "
$ws.Range("D6").Value = "'167.09"
$ws.Range("E6").Value = "  -1.55%  "
$ws.Range("D7").Value = "3.837.03"
$ws.Range("E7").Value = "  -2.51%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("E9").Value = "  -1.51%  "
$ws.Range("E10").Value = "  -5.27%  "
$ws.Range("E11").Value = "  -0.77%  "
$ws.Range("E12").Value = "  -2.77%  "
$ws.Range("D13").Value = "'0.0000259"
$ws.Range("E13").Value = "  +0.60%  "
$ws.Range("D14").Value = "'36.81"
$ws.Range("E14").Value = "  -3.20%  "
$ws.Range("D15").Value = "4.479.00"
$ws.Range("E15").Value = "  -2.47%  "
$ws.Range("D16").Value = "3.826.25"
$ws.Range("E16").Value = "  -2.50%  "
$ws.Range("D17").Value = "68.095.70"
$ws.Range("E17").Value = "  -2.41%  "
$ws.Range("D18").Value = "'17.95"
$ws.Range("D19").Value = "'7.31"
$ws.Range("E19").Value = "  -3.53%  "
$ws.Range("E20").Value = "  -0.59%  "
$ws.Range("E21").Value = "  -3.69%  "
$ws.Range("D22").Value = "'464.91"
$ws.Range("E22").Value = "  -6.99%  "
$ws.Range("E23").Value = "  -1.24%  "
$ws.Range("D24").Value = "'0.0000159"
$ws.Range("E24").Value = "  -6.32%  "
$ws.Range("D25").Value = "'82.76"
$ws.Range("E25").Value = "  -3.20%  "
$ws.Range("D26").Value = "'2.22"
$ws.Range("E26").Value = "  -3.36%  "
$ws.Range("D27").Value = "'12.00"
$ws.Range("E27").Value = "  -2.78%  "
$ws.Range("E28").Value = "  -0.22%  "
$ws.Range("D29").Value = "'9.95"
$ws.Range("E29").Value = "  -3.28%  "
$ws.Range("D30").Value = "'2.96"
$ws.Range("E30").Value = "  -1.30%  "
$ws.Range("D31").Value = "3.985.41"
$ws.Range("E31").Value = "  -2.44%  "
$ws.Range("D32").Value = "'7.63"
$ws.Range("E32").Value = "  -3.03%  "
$ws.Range("E33").Value = "  -5.05%  "
$ws.Range("D34").Value = "'31.08"
$ws.Range("E34").Value = "  -3.72%  "
$ws.Range("D35").Value = "'9.38"
$ws.Range("E35").Value = "  -2.33%  "
$ws.Range("D36").Value = "3.802.33"
$ws.Range("E36").Value = "  -2.62%  "
$ws.Range("E37").Value = "  -3.73%  "
$ws.Range("D38").Value = "'3.60"
$ws.Range("E38").Value = "  +8.56%  "
$ws.Range("E39").Value = "  -0.23%  "
$ws.Range("E40").Value = "  -2.13%  "
$ws.Range("D41").Value = "'5.88"
$ws.Range("E41").Value = "  -4.39%  "
$ws.Range("E42").Value = "  +0.27%  "
$ws.Range("D43").Value = "'0.312"
$ws.Range("E43").Value = "  -4.13%  "
$ws.Range("D44").Value = "'420.83"
$ws.Range("E44").Value = "  -3.55%  "
$ws.Range("E45").Value = "  -4.95%  "
$ws.Range("E47").Value = "  +4.37%  "
$ws.Range("B48").Value = "OKB"
$ws.Range("C48").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D48").Value = "'47.03"
$ws.Range("E48").Value = "  -2.61%  "
$ws.Range("B49").Value = "Cosmos"
$ws.Range("C49").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D49").Value = "'8.56"
$ws.Range("E49").Value = "  -0.61%  "
$ws.Range("D50").Value = "'142.72"
$ws.Range("E50").Value = "  -0.33%  "
$ws.Range("D51").Value = "'26.02"
$ws.Range("E51").Value = "  +0.60%  "
